Write-Output $ppt
Write-Output $ppt.Name
Write-Output $app
